$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Each entry: cell reference, new value, and whether the value must be
# forced to Text (numeric-looking strings would otherwise be parsed as numbers).
$changes = @(
    @('D2', '61.458.74', 0),
    @('E2', '  -1.43%  ', 0),
    @('D3', '3.351.65', 0),
    @('E3', '  -0.46%  ', 0),
    @('E4', '  -0.10%  ', 0),
    @('D5', '400.44', 1),
    @('E5', '  -3.46%  ', 0),
    @('D6', '125.97', 1),
    @('E6', '  +7.94%  ', 0),
    @('D7', '0.588', 1),
    @('E7', '  +1.96%  ', 0),
    @('D8', '1.00', 1),
    @('E8', '  +0.02%  ', 0),
    @('D9', '0.659', 1),
    @('E9', '  +4.51%  ', 0),
    @('D10', '0.119', 1),
    @('E10', '  +1.27%  ', 0),
    @('D11', '40.77', 1),
    @('E11', '  +1.63%  ', 0),
    @('E12', '  -0.98%  ', 0),
    @('D13', '3.874.44', 0),
    @('E13', '  -0.29%  ', 0),
    @('D14', '8.26', 1),
    @('E14', '  -0.92%  ', 0),
    @('D15', '19.25', 1),
    @('E15', '  -0.43%  ', 0),
    @('D16', '3.334.74', 0),
    @('E16', '  -0.60%  ', 0),
    @('D17', '61.362.77', 0),
    @('E17', '  -1.22%  ', 0),
    @('D18', '11.20', 1),
    @('E18', '  +2.99%  ', 0),
    @('E19', '  -0.63%  ', 0),
    @('D20', '0.0000127', 1),
    @('E20', '  +7.62%  ', 0),
    @('E21', '  -4.53%  ', 0),
    @('D22', '79.80', 1),
    @('E22', '  +6.62%  ', 0),
    @('E23', '  +0.89%  ', 0),
    @('D24', '298.55', 1),
    @('E24', '  +0.83%  ', 0),
    @('D25', '3.10', 1),
    @('E25', '  -1.29%  ', 0),
    @('E26', '  +11.22%  ', 0),
    @('B27', 'Filecoin', 0),
    @('C27', 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil', 0),
    @('D27', '8.19', 1),
    @('E27', '  +7.58%  ', 0),
    @('B28', 'EthereumClassic', 0),
    @('C28', 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc', 0),
    @('D28', '28.94', 1),
    @('E28', '  -1.58%  ', 0),
    @('D29', '7.43', 1),
    @('E29', '  -6.80%  ', 0),
    @('E30', '  -2.03%  ', 0),
    @('D31', '0.114', 1),
    @('E31', '  +0.22%  ', 0),
    @('E32', '  +0.09%  ', 0),
    @('D33', '11.31', 1),
    @('E33', '  -1.33%  ', 0),
    @('E34', '  -2.07%  ', 0),
    @('D35', '40.96', 1),
    @('E35', '  -5.13%  ', 0),
    @('D36', '0.0477', 1),
    @('E36', '  -2.94%  ', 0),
    @('D37', '51.96', 1),
    @('E37', '  -0.63%  ', 0),
    @('E38', '  +0.15%  ', 0),
    @('D39', '3.36', 1),
    @('E39', '  -2.59%  ', 0),
    @('D40', '2.90', 1),
    @('E40', '  -7.17%  ', 0),
    @('D41', '136.93', 1),
    @('E41', '  +2.68%  ', 0),
    @('D42', '1.95', 1),
    @('E42', '  +2.30%  ', 0),
    @('E43', '  +1.03%  ', 0),
    @('D44', '0.282', 1),
    @('E44', '  -1.81%  ', 0),
    @('D45', '3.88', 1),
    @('E45', '  +0.13%  ', 0),
    @('D46', '16.53', 1),
    @('E46', '  +0.20%  ', 0),
    @('D47', '2.22', 1),
    @('E47', '  -0.41%  ', 0),
    @('D48', '20.99', 1),
    @('E48', '  -1.07%  ', 0),
    @('D49', '3.677.67', 0),
    @('E49', '  -0.14%  ', 0),
    @('D50', '2.092.04', 0),
    @('E50', '  -3.52%  ', 0),
    @('E51', '  -4.85%  ', 0),
)

foreach ($change in $changes) {
    $cellRef = $change[0]
    $value = $change[1]
    $forceText = $change[2]
    $range = $ws.Range($cellRef)
    if ($forceText -eq 1) {
        # Prefix with an apostrophe so Excel stores the numeric-looking
        # string as text instead of converting it to a number, then
        # restore the default "Normal" style so no formatting side
        # effect (e.g. a Text number format) is left behind.
        $range.Value = "'" + $value
        $range.Style = "Normal"
    } else {
        $range.Value = $value
    }
}
